$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.665.86"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "2.894.86"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'566.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").Value = "'143.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "2.892.20"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").Value = "'6.95"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").Value = "'0.429"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "'31.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "'0.125"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "3.375.50"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "61.652.63"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "'6.52"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "2.895.27"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "'430.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "'13.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").Value = "'0.653"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'79.14"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'11.85"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -11.29%  "
$ws.Range("E28").Value = "  -5.63%  "
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("D30").Value = "'7.03"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("E31").Value = "  -4.31%  "
$ws.Range("D32").Value = "'2.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.74%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'25.51"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("D36").Value = "'0.960"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("D37").Value = "'5.36"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -5.10%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.55%  "
$ws.Range("D41").Value = "'8.18"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("E42").Value = "  -4.06%  "
$ws.Range("D43").Value = "'39.81"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  -4.23%  "
$ws.Range("D45").Value = "2.688.46"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "'132.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "'343.20"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.75%  "
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "'21.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.40%  "
